$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 14 from 45185 to 45204
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
